$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1913.0233
$ws.Range("I137").Value = 1947.037
$ws.Range("J137").Value = 1855.625
$ws.Range("K137").Value = 5841.111
$ws.Range("L137").Value = 5566.875
$ws.Range("M137").Value = -3291.111
$ws.Range("N137").Value = -10666.875
$ws.Range("H141").Value = 22200.715
$ws.Range("I141").Value = 5100
$ws.Range("J141").Value = 35026.25
$ws.Range("K141").Value = 15300
$ws.Range("L141").Value = 105078.75
$ws.Range("M141").Value = -10120
$ws.Range("N141").Value = -115438.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 50000
$ws.Range("J34").Value = 50000
$ws.Range("L34").Value = 50000
$ws.Range("N34").Value = -50542
$ws.Range("H61").Value = 3986.4666
$ws.Range("I61").Value = 2096.5557
$ws.Range("J61").Value = 6821.3335
$ws.Range("K61").Value = 2096.5557
$ws.Range("L61").Value = 6821.3335
$ws.Range("M61").Value = -1884.5557
$ws.Range("N61").Value = -7245.3335
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("H122").Value = 1542.4
$ws.Range("I122").Value = 1542.4
$ws.Range("K122").Value = 4627.200000000001
$ws.Range("M122").Value = -2177.200000000001
$ws.Range("H132").Value = 4096.1523
$ws.Range("I132").Value = 1536.2812
$ws.Range("J132").Value = 9947.286
$ws.Range("K132").Value = 4608.8436
$ws.Range("L132").Value = 29841.858
$ws.Range("M132").Value = -2078.8436
$ws.Range("N132").Value = -34901.858
$ws.Range("H136").Value = 3986.4666
$ws.Range("I136").Value = 2096.5557
$ws.Range("J136").Value = 6821.3335
$ws.Range("K136").Value = 6289.6671
$ws.Range("L136").Value = 20464.0005
$ws.Range("M136").Value = -3739.6671
$ws.Range("N136").Value = -25564.0005
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 234190.8
$ws.Range("J47").Value = 234190.8
$ws.Range("L47").Value = 234190.8
$ws.Range("N47").Value = -235230.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 40000
$ws.Range("J18").Value = 40000
$ws.Range("L18").Value = 40000
$ws.Range("N18").Value = -40460
$ws.Range("H31").Value = 6412293
$ws.Range("I31").Value = 1539.881
$ws.Range("J31").Value = 33337456
$ws.Range("K31").Value = 1539.881
$ws.Range("L31").Value = 33337456
$ws.Range("M31").Value = -1244.881
$ws.Range("N31").Value = -33338046
$ws.Range("H34").Value = 6412293
$ws.Range("I34").Value = 1539.881
$ws.Range("J34").Value = 33337456
$ws.Range("K34").Value = 1539.881
$ws.Range("L34").Value = 33337456
$ws.Range("M34").Value = -1337.881
$ws.Range("N34").Value = -33337860
$ws.Range("H58").Value = 1319104.5
$ws.Range("I58").Value = 2509.7083
$ws.Range("K58").Value = 2509.7083
$ws.Range("M58").Value = -2306.7083
$ws.Range("H132").Value = 3705.5557
$ws.Range("I132").Value = 2187.111
$ws.Range("J132").Value = 5224
$ws.Range("K132").Value = 6561.333
$ws.Range("L132").Value = 15672
$ws.Range("M132").Value = -4031.333
$ws.Range("N132").Value = -20732
$ws.Range("H134").Value = 2654
$ws.Range("I134").Value = 1278.4
$ws.Range("J134").Value = 4946.6665
$ws.Range("K134").Value = 3835.2
$ws.Range("L134").Value = 14839.9995
$ws.Range("M134").Value = -1300.2
$ws.Range("N134").Value = -19909.9995
$ws.Range("H136").Value = 1319104.5
$ws.Range("I136").Value = 2509.7083
$ws.Range("K136").Value = 7529.124899999999
$ws.Range("M136").Value = -4979.124899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3000
$ws.Range("J110").Value = 3960
$ws.Range("L110").Value = 11880
$ws.Range("N110").Value = -20060
$ws.Range("H120").Value = 38471624
$ws.Range("I120").Value = 83337350
$ws.Range("J120").Value = 15290
$ws.Range("K120").Value = 250012050
$ws.Range("L120").Value = 45870
$ws.Range("M120").Value = -250007212
$ws.Range("N120").Value = -55546

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2846.3125
$ws.Range("I80").Value = 2636.25
$ws.Range("J80").Value = 3476.5
$ws.Range("K80").Value = 2636.25
$ws.Range("L80").Value = 3476.5
$ws.Range("M80").Value = -1638.25
$ws.Range("N80").Value = -5472.5
$ws.Range("H83").Value = 2846.3125
$ws.Range("I83").Value = 2636.25
$ws.Range("J83").Value = 3476.5
$ws.Range("K83").Value = 13181.25
$ws.Range("L83").Value = 17382.5
$ws.Range("M83").Value = -8189.25
$ws.Range("N83").Value = -27366.5
$ws.Range("H108").Value = 39995
$ws.Range("J108").Value = 39995
$ws.Range("L108").Value = 39995
$ws.Range("N108").Value = -47675

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 13500
$ws.Range("I88").Value = 8000
$ws.Range("J88").Value = 19000
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 19000
$ws.Range("M88").Value = -7572
$ws.Range("N88").Value = -19856
$ws.Range("H91").Value = 13500
$ws.Range("I91").Value = 8000
$ws.Range("J91").Value = 19000
$ws.Range("K91").Value = 8000
$ws.Range("L91").Value = 19000
$ws.Range("M91").Value = -6518
$ws.Range("N91").Value = -21964

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 25000
$ws.Range("J16").Value = 25000
$ws.Range("L16").Value = 25000
$ws.Range("N16").Value = -25584
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H109").Value = 23289
$ws.Range("J109").Value = 23289
$ws.Range("L109").Value = 23289
$ws.Range("N109").Value = -26063
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H136").Value = 6581069.5
$ws.Range("I136").Value = 13159444
$ws.Range("J136").Value = 2694.7368
$ws.Range("K136").Value = 39478332
$ws.Range("L136").Value = 8084.2104
$ws.Range("M136").Value = -39475782
$ws.Range("N136").Value = -13184.2104
$ws.Range("N108").ClearContents()
$ws.Range("N110").ClearContents()
$ws.Range("N133").ClearContents()
